# Update assignment due dates:
#   - "COS318 - FA2018" header stays the same text, but the stray
#     "_GoBack" bookmark that used to sit at the end of that paragraph
#     is removed ...
#   - ... and re-created after the due date is updated, right after the
#     new year digit (mirrors a user editing the due-date line last).
#   - "October 19th, 2017" -> "October 18th, 2018"

$d = $word.ActiveDocument

# --- Remove the old _GoBack bookmark (end of the "FA2018" paragraph) ---
$d.Bookmarks.Item("_GoBack").Delete()

# --- "19th" -> "18th" -------------------------------------------------
# Change only the "9" so the surrounding text/runs (incl. the
# superscripted "th") are left alone. Toggling a character property
# around the text assignment keeps the edited digit in its own run
# instead of it being silently re-merged into its neighbour.
$full = $d.Content.Text
$idx = $full.IndexOf("19th")
$rNine = $d.Range($idx + 1, $idx + 2)
$rNine.Bold = 1
$rNine.Text = "8"
$full = $d.Content.Text
$idx = $full.IndexOf("18th")
$rEight = $d.Range($idx + 1, $idx + 2)
$rEight.Bold = 0

# --- ", 2017" -> ", 2018" ----------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf(", 2017")
$rSeven = $d.Range($idx + 5, $idx + 6)
$rSeven.Bold = 1
$rSeven.Text = "8"
$full = $d.Content.Text
$idx = $full.IndexOf(", 2018")
$rEight2 = $d.Range($idx + 5, $idx + 6)
$rEight2.Bold = 0

# --- Re-add the _GoBack bookmark right after the new "2018" year, -----
# --- before the line break that follows it -----------------------------
$full = $d.Content.Text
$idx = $full.IndexOf(", 2018")
$rPoint = $d.Range($idx + 6, $idx + 6)
$d.Bookmarks.Add("_GoBack", $rPoint)

Write-Output $d.Paragraphs(1).Range.Text
Write-Output $d.Paragraphs(2).Range.Text
